$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = @(-8.5100471597281082,-7.9086441958729354,-10.329385548761364,-17.091056618930654,-10.397201348853935,-7.8172393134921947,-8.6243083469668633,-9.4245675938825464,-9.673843579769807,-8.5111787528298084,-15.579798565294226,-8.4983189785559947,-12.149809226904312,-9.0603207717813081,-11.708496338834705,-9.8385621332702335,-8.598560757742769,-9.9303635603070095,-10.998509661049008,-7.7154093207080301,-7.3462527672627411,-7.3537937885633484,-8.7602543021568415,-8.168055548206544,-9.2515621906550383,-7.759137069280265,-8.4967647042250558,-12.83313590912852,-7.8627114935513793,-13.048480706831526,-10.721129966845455,-7.1768145875864207,-11.989187220830736,-8.8972697157594389,-9.5381757723288541,-10.913175528997597,-9.1912983543465359,-10.075051442458967,-9.5352229544555893,-9.6086518444576239,-9.5750032368753928,-8.6724126136157444,-9.8978472111921754,-12.051269252788147,-8.0854126856777171,-8.857343182423957,-8.9811906150788516,-7.6334034856355109,-12.980307113231895,-8.8984118378081103,-12.175221925075027,-9.9531389055999817,-8.6642210334633081,-10.183295859526689,-13.645990853995054,-11.446584147825897,-9.191506004892684,-9.3701111321284341,-8.0291240314436081,-13.924557401088741,-7.1491876397947545,-8.4424872545712226,-10.991807936825781,-9.0912043317857041,-12.260629327869809,-7.952004009199821,-8.3688704033604555,-11.166272750596564,-8.2074545209074774,-10.683596075951183,-8.0553113402317411,-7.444254796887722,-8.1983667931254622,-10.299932791335669,-9.3976171131063335,-9.8184832800875235,-8.6796517947407121,-8.5728059691948442,-13.781053359743893,-8.0830948759847931,-7.8716564998910465,-17.131463231895573,-7.7532951856315373,-6.8820847484721046,-8.3743206699066217,-8.7981773369513991,-13.767845771279294,-14.944642547299559,-9.1147850353919573,-10.675858094392632,-9.6572191381899355,-8.3139855586730249,-12.054807853849072,-8.8516120793253812,-8.083777604752413,-9.7087973557063236,-8.4862656260277909,-8.1905889893343335,-12.371168587638099,-9.2684529660664161)
$row2 = @(-7.0592767533815826,-8.3778464743528307,-9.7648928169478726,-15.214022140843413,-8.8489107302980266,-7.3078539094174646,-8.1856710245864832,-8.9360681516936555,-9.1417521560362029,-8.0363125101067165,-14.803191060683458,-8.9718000359819037,-10.528942197597383,-8.5140764151411794,-11.071623745698718,-9.3908715473508533,-8.1692170984846069,-9.3849769028600285,-10.518054116908049,-6.32482798056368,-7.89252419001195,-7.8681740343909867,-8.2564726383375451,-7.7106796317149016,-8.745616337077891,-7.3456045197961393,-8.0566860417786135,-12.182751032669701,-8.2951354006425362,-12.406321775424304,-10.165284490486043,-7.6918216350264821,-10.350850217501609,-8.3977151009099007,-9.0409901088586295,-9.2879455057777829,-8.6825550129277111,-9.5412033016199427,-8.9985200158126872,-9.1172802949527405,-9.0428955797773654,-8.1221066939401556,-10.344892618770055,-11.364425180279623,-7.6445492162375883,-8.3300338756488959,-8.4624124816694621,-8.1970101972560663,-12.273521202224311,-9.3379301935953372,-10.579665766124711,-8.4272047557276473,-8.1563543238522698,-9.6397960283320501,-12.942011879768312,-10.80904810737403,-8.7362875662206232,-8.8594314742738742,-7.5833775008687008,-13.280304148576247,-7.7647653588368444,-7.9620931472853105,-10.340331544488253,-7.5797221029201109,-11.573079459796393,-7.5382671702738957,-7.8496925892224363,-9.6343803619671498,-7.7503712216127099,-10.152198379362636,-7.5924010921037235,-8.0272927058129859,-6.7001675793182747,-8.7703703548526395,-9.7818855455673184,-9.2439869308678784,-8.2972922752613183,-8.088909070795717,-12.106300160133385,-7.6168611584898045,-7.4476188353076909,-16.253318045852485,-7.3220913198605224,-7.4412763538931053,-7.9125793609540676,-8.2613459162564524,-12.11592865212212,-14.251122478085195,-8.6376546001739545,-9.1949044082248381,-9.1124608559650326,-7.8089654855984074,-11.329935117898764,-8.330361227248293,-7.5804771639370596,-9.1721227834123997,-8.0310726656999609,-7.728530600027141,-11.744999259078279,-8.6895250762246636)

for ($i = 0; $i -lt $row1.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $row1[$i]
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

Write-Host "Updated $($row1.Length) columns for rows 1 and 2"